# Generate Report for Handoff
#
# The localization run id (GUID) changes from 35314411-9543-4f97-b74a-531f3dc4044d
# to 853f99db-029a-4e31-8bac-f1d20846b513, the handoff package hash changes from
# ac4e95baadf25df9c159eac89c4f5645e8135645 to 5fcd11e6ec0671009cc70f18cfae2e578915ca1a,
# and the handoff timestamps advance by ~30s on each of the two language sheets.
# This touches the source-markdown / handoff-xlf file names shown (and hyperlinked)
# on all three sheets, plus the "Latest Handoff Datetime" cell on the two language
# sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "35314411-9543-4f97-b74a-531f3dc4044d"
$newGuid = "853f99db-029a-4e31-8bac-f1d20846b513"
$oldHash = "ac4e95baadf25df9c159eac89c4f5645e8135645"
$newHash = "5fcd11e6ec0671009cc70f18cfae2e578915ca1a"

$oldMdName  = "$oldGuid.md"
$newMdName  = "$newGuid.md"
$oldZhName  = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhName  = "$newGuid.$newHash.zh-cn.xlf"
$oldDeName  = "$oldGuid.$oldHash.de-de.xlf"
$newDeName  = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": A2 links to the source .md file, A3 to .localization-config
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$mdAddrOverview     = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/e2e/$oldMdName"
$configAddrOverview = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/.localization-config"

# Rewrite every hyperlink on the sheet (re-adding with the SAME underlying
# target address, just the new display text) so the "display" attribute
# tracks the renamed file instead of going stale.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddrOverview, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddrOverview, "", "", ".localization-config")

$wsOverview.Range("A2").Value2 = $newMdName

# ---------------------------------------------------------------------------
# Sheet "zh-cn": A2 -> source .md, C2 -> handoff .xlf, D2 -> handoff datetime
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$mdAddrZhCn     = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/e2e/$oldMdName"
$xlfAddrZhCn    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9807bc5eee4381db80cfdb66dbcc82ce938077b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhName"
$configAddrZhCn = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/.localization-config"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddrZhCn, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $xlfAddrZhCn, "", "", $newZhName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configAddrZhCn, "", "", ".localization-config")

$wsZhCn.Range("A2").Value2 = $newMdName
$wsZhCn.Range("C2").Value2 = $newZhName
$wsZhCn.Range("D2").Value2 = "2016-03-09 14:26:00"

# ---------------------------------------------------------------------------
# Sheet "de-de": A2 -> source .md, C2 -> handoff .xlf, D2 -> handoff datetime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$mdAddrDeDe     = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/e2e/$oldMdName"
$xlfAddrDeDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/756c5baa62b2312f40399dce92e3923d1c29f097/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeName"
$configAddrDeDe = "https://github.com/OpenLocalizationTest/oltest/blob/683c17e2da3924c66e595dcae48852c9fd4b7a52/.localization-config"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddrDeDe, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $xlfAddrDeDe, "", "", $newDeName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configAddrDeDe, "", "", ".localization-config")

$wsDeDe.Range("A2").Value2 = $newMdName
$wsDeDe.Range("C2").Value2 = $newDeName
$wsDeDe.Range("D2").Value2 = "2016-03-09 14:26:05"
